$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "30.222.14"
Set-TextValue $ws.Range("E2") "  -0.59%  "
Set-TextValue $ws.Range("D3") "2.073.81"
Set-TextValue $ws.Range("E3") "  +3.17%  "
Set-TextValue $ws.Range("D4") "0.9992"
Set-TextValue $ws.Range("E4") "  -0.23%  "
Set-TextValue $ws.Range("D5") "327.14"
Set-TextValue $ws.Range("E5") "  +0.70%  "
Set-TextValue $ws.Range("E6") "  -0.20%  "
Set-TextValue $ws.Range("D7") "0.5187"
Set-TextValue $ws.Range("E7") "  +1.60%  "
Set-TextValue $ws.Range("D8") "0.4321"
Set-TextValue $ws.Range("E8") "  +3.82%  "
Set-TextValue $ws.Range("D9") "0.08839"
Set-TextValue $ws.Range("E9") "  +0.60%  "
Set-TextValue $ws.Range("D10") "45.69"
Set-TextValue $ws.Range("E10") "  +6.43%  "
Set-TextValue $ws.Range("D12") "24.18"
Set-TextValue $ws.Range("E12") "  -2.07%  "
Set-TextValue $ws.Range("D13") "2.070.44"
Set-TextValue $ws.Range("E13") "  +3.37%  "
Set-TextValue $ws.Range("D14") "6.649"
Set-TextValue $ws.Range("E14") "  +0.73%  "
Set-TextValue $ws.Range("D15") "7.656"
Set-TextValue $ws.Range("E15") "  +2.31%  "
Set-TextValue $ws.Range("D16") "94.88"
Set-TextValue $ws.Range("E16") "  +0.58%  "
Set-TextValue $ws.Range("D17") "0.9995"
Set-TextValue $ws.Range("E17") "  -0.31%  "
Set-TextValue $ws.Range("D18") "0.00001119"
Set-TextValue $ws.Range("E18") "  +0.11%  "
Set-TextValue $ws.Range("E19") "  +0.99%  "
Set-TextValue $ws.Range("D20") "18.72"
Set-TextValue $ws.Range("E20") "  -1.42%  "
Set-TextValue $ws.Range("D21") "0.9983"
Set-TextValue $ws.Range("E21") "  -0.17%  "
Set-TextValue $ws.Range("D22") "6.203"
Set-TextValue $ws.Range("E22") "  -0.40%  "
Set-TextValue $ws.Range("D23") "30.269.33"
Set-TextValue $ws.Range("E23") "  -0.61%  "
Set-TextValue $ws.Range("D24") "12.26"
Set-TextValue $ws.Range("E24") "  +2.29%  "
Set-TextValue $ws.Range("D25") "2.281"
Set-TextValue $ws.Range("E25") "  +2.41%  "
Set-TextValue $ws.Range("D26") "2.312.64"
Set-TextValue $ws.Range("E26") "  +3.40%  "
Set-TextValue $ws.Range("D27") "22.14"
Set-TextValue $ws.Range("E27") "  -1.00%  "
Set-TextValue $ws.Range("D28") "2.535"
Set-TextValue $ws.Range("E28") "  +4.84%  "
Set-TextValue $ws.Range("D29") "161.44"
Set-TextValue $ws.Range("E29") "  -0.99%  "
Set-TextValue $ws.Range("D30") "130.53"
Set-TextValue $ws.Range("E30") "  -0.76%  "
Set-TextValue $ws.Range("E31") "  +4.37%  "
Set-TextValue $ws.Range("E32") "  +1.03%  "
Set-TextValue $ws.Range("B33") "ARBITRUM"
Set-TextValue $ws.Range("C33") "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue $ws.Range("D33") "1.573"
Set-TextValue $ws.Range("E33") "  +16.41%  "
Set-TextValue $ws.Range("B34") "Filecoin"
Set-TextValue $ws.Range("C34") "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue $ws.Range("D34") "6.074"
Set-TextValue $ws.Range("E34") "  -0.78%  "
Set-TextValue $ws.Range("D35") "3.830"
Set-TextValue $ws.Range("E35") "  -0.05%  "
Set-TextValue $ws.Range("D36") "0.02561"
Set-TextValue $ws.Range("E36") "  +1.49%  "
Set-TextValue $ws.Range("D37") "9.633"
Set-TextValue $ws.Range("E37") "  +5.53%  "
Set-TextValue $ws.Range("B38") "InternetComputer(DFINITY)"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextValue $ws.Range("D38") "5.400"
Set-TextValue $ws.Range("E38") "  -1.14%  "
Set-TextValue $ws.Range("B39") "Hedera"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.06599"
Set-TextValue $ws.Range("E39") "  -0.22%  "
Set-TextValue $ws.Range("D40") "12.53"
Set-TextValue $ws.Range("E40") "  +0.75%  "
Set-TextValue $ws.Range("D41") "0.2234"
Set-TextValue $ws.Range("E41") "  +1.54%  "
Set-TextValue $ws.Range("D42") "0.6754"
Set-TextValue $ws.Range("E42") "  +1.24%  "
Set-TextValue $ws.Range("D43") "1.242"
Set-TextValue $ws.Range("E43") "  +0.81%  "
Set-TextValue $ws.Range("D44") "0.9978"
Set-TextValue $ws.Range("E44") "  -0.22%  "
Set-TextValue $ws.Range("B45") "Decentraland"
Set-TextValue $ws.Range("C45") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D45") "0.6335"
Set-TextValue $ws.Range("E45") "  +2.51%  "
Set-TextValue $ws.Range("B46") "EnergySwap"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D46") "13.83"
Set-TextValue $ws.Range("E46") "  +1.50%  "
Set-TextValue $ws.Range("D47") "2.192"
Set-TextValue $ws.Range("E47") "  -0.26%  "
Set-TextValue $ws.Range("D48") "3.597"
Set-TextValue $ws.Range("E48") "  -1.95%  "
Set-TextValue $ws.Range("D49") "1.233"
Set-TextValue $ws.Range("E49") "  -2.87%  "
Set-TextValue $ws.Range("D50") "1.182"
Set-TextValue $ws.Range("E50") "  +6.64%  "
Set-TextValue $ws.Range("D51") "81.16"
Set-TextValue $ws.Range("E51") "  +0.05%  "
